$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.37
$ws.Range("D4").Value = -7.795999999999999
$ws.Range("C6").Value = -11.714
$ws.Range("C7").Value = -12.995
$ws.Range("D9").Value = -8.026999999999999
$ws.Range("D12").Value = -7.451000000000001
$ws.Range("C16").Value = -12.429
$ws.Range("D17").Value = -8.464
$ws.Range("D18").Value = -8.516999999999999
$ws.Range("D19").Value = -7.754
$ws.Range("C20").Value = -12.182
$ws.Range("D20").Value = -7.375
$ws.Range("D26").Value = -7.602000000000001
$ws.Range("C28").Value = -12.409
$ws.Range("C29").Value = -12.148
$ws.Range("D31").Value = -7.991
$ws.Range("C32").Value = -12.475
$ws.Range("D39").Value = -7.617999999999999
$ws.Range("C40").Value = -11.745
$ws.Range("D40").Value = -7.558
$ws.Range("D41").Value = -7.928
$ws.Range("D42").Value = -8.036
$ws.Range("D43").Value = -7.946999999999998
$ws.Range("C46").Value = -13.421
$ws.Range("D47").Value = -7.475
$ws.Range("D48").Value = -7.584999999999999
$ws.Range("C51").Value = -11.34
$ws.Range("C52").Value = -11.218
$ws.Range("C57").Value = -13.848
$ws.Range("C59").Value = -11.711
$ws.Range("C62").Value = -13.514
$ws.Range("D63").Value = -7.25
$ws.Range("D64").Value = -7.580999999999999
$ws.Range("C66").Value = -11.713
$ws.Range("C73").Value = -12.573
$ws.Range("C74").Value = -12.078
$ws.Range("D76").Value = -7.637
$ws.Range("D81").Value = -8.117000000000001
$ws.Range("D89").Value = -7.986999999999999
$ws.Range("C92").Value = -10.49
$ws.Range("D94").Value = -7.593999999999999
$ws.Range("C100").Value = -11.221
